$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 237. This pushes the existing rows 237-271
# down to 238-272, preserving all of their data/formatting untouched.
$ws.Rows("237:237").Insert()

# Populate the new row 237 with its own data (a new price-record entry).
$ws.Cells.Item(237, 1).Value = 4
$ws.Cells.Item(237, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(237, 3).Value = "Los Lagos"
$ws.Cells.Item(237, 4).Value = 44491
$ws.Cells.Item(237, 5).Value = 10
$ws.Cells.Item(237, 6).Value = 100112006
$ws.Cells.Item(237, 7).Value = "Repollo"
$ws.Cells.Item(237, 8).Value = "Crespo record"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 1500
$ws.Cells.Item(237, 11).Value = 1200
$ws.Cells.Item(237, 12).Value = 1200
$ws.Cells.Item(237, 13).Value = 1200
$ws.Cells.Item(237, 14).Value = "$/unidad"
$ws.Cells.Item(237, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(237, 16).Value = 1200
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = "Hortaliza"
